$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table grows from 8 data rows (sending clusters FAPs/MuSCs x target FAPs/MuSCs/ECs/Resolving-Mac)
# to 9 data rows (sending clusters ECs/FAPs/MuSCs x target ECs/FAPs/MuSCs), with refreshed TPM-based values.
# Clear the old 2:9 block first, then write the new 2:10 block.
$ws.Range("A2:T9").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Edn3"
$ws.Range("C2").Value = "Ednra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1264883333333333
$ws.Range("H2").Value = 0.379465
$ws.Range("I2").Value = 0.02088586470611676
$ws.Range("J2").Value = 0.02088586470611676
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.169998666666667
$ws.Range("N2").Value = 9.509996000000001
$ws.Range("O2").Value = 0.06457634599094531
$ws.Range("P2").Value = 0.06457634599094531
$ws.Range("Q2").Value = 0.4009678480155556
$ws.Range("R2").Value = 3.60871063214
$ws.Range("S2").Value = 0.001348732825582269
$ws.Range("T2").Value = 0.001348732825582269

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Edn3"
$ws.Range("C3").Value = "Ednra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1264883333333333
$ws.Range("H3").Value = 0.379465
$ws.Range("I3").Value = 0.02088586470611676
$ws.Range("J3").Value = 0.02088586470611676
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 32.709374
$ws.Range("N3").Value = 98.12812199999999
$ws.Range("O3").Value = 0.6663257858061865
$ws.Range("P3").Value = 0.6663257858061865
$ws.Range("Q3").Value = 4.137354201636667
$ws.Range("R3").Value = 37.23618781472999
$ws.Range("S3").Value = 0.01391679021254495
$ws.Range("T3").Value = 0.01391679021254495

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Edn3"
$ws.Range("C4").Value = "Ednra"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1264883333333333
$ws.Range("H4").Value = 0.379465
$ws.Range("I4").Value = 0.02088586470611676
$ws.Range("J4").Value = 0.02088586470611676
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.20978866666667
$ws.Range("N4").Value = 39.629366
$ws.Range("O4").Value = 0.2690978682028682
$ws.Range("P4").Value = 0.2690978682028682
$ws.Range("Q4").Value = 1.670884152132222
$ws.Range("R4").Value = 15.03795736919
$ws.Range("S4").Value = 0.005620341667989544
$ws.Range("T4").Value = 0.005620341667989544

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Edn3"
$ws.Range("C5").Value = "Ednra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.012975
$ws.Range("H5").Value = 0.038925
$ws.Range("I5").Value = 0.002142443397113291
$ws.Range("J5").Value = 0.002142443397113291
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.169998666666667
$ws.Range("N5").Value = 9.509996000000001
$ws.Range("O5").Value = 0.06457634599094531
$ws.Range("P5").Value = 0.06457634599094531
$ws.Range("Q5").Value = 0.04113073270000001
$ws.Range("R5").Value = 0.3701765943000001
$ws.Range("S5").Value = 0.0001383511660780041
$ws.Range("T5").Value = 0.0001383511660780041

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Edn3"
$ws.Range("C6").Value = "Ednra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.012975
$ws.Range("H6").Value = 0.038925
$ws.Range("I6").Value = 0.002142443397113291
$ws.Range("J6").Value = 0.002142443397113291
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 32.709374
$ws.Range("N6").Value = 98.12812199999999
$ws.Range("O6").Value = 0.6663257858061865
$ws.Range("P6").Value = 0.6663257858061865
$ws.Range("Q6").Value = 0.42440412765
$ws.Range("R6").Value = 3.81963714885
$ws.Range("S6").Value = 0.001427565280126789
$ws.Range("T6").Value = 0.001427565280126789

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Edn3"
$ws.Range("C7").Value = "Ednra"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.012975
$ws.Range("H7").Value = 0.038925
$ws.Range("I7").Value = 0.002142443397113291
$ws.Range("J7").Value = 0.002142443397113291
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.20978866666667
$ws.Range("N7").Value = 39.629366
$ws.Range("O7").Value = 0.2690978682028682
$ws.Range("P7").Value = 0.2690978682028682
$ws.Range("Q7").Value = 0.17139700795
$ws.Range("R7").Value = 1.54257307155
$ws.Range("S7").Value = 0.0005765269509084974
$ws.Range("T7").Value = 0.0005765269509084975

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Edn3"
$ws.Range("C8").Value = "Ednra"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.916706
$ws.Range("H8").Value = 17.750118
$ws.Range("I8").Value = 0.9769716918967699
$ws.Range("J8").Value = 0.97697169189677
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.169998666666667
$ws.Range("N8").Value = 9.509996000000001
$ws.Range("O8").Value = 0.06457634599094531
$ws.Range("P8").Value = 0.06457634599094531
$ws.Range("Q8").Value = 18.75595013105867
$ws.Range("R8").Value = 168.803551179528
$ws.Range("S8").Value = 0.06308926199928502
$ws.Range("T8").Value = 0.06308926199928504

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Edn3"
$ws.Range("C9").Value = "Ednra"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.916706
$ws.Range("H9").Value = 17.750118
$ws.Range("I9").Value = 0.9769716918967699
$ws.Range("J9").Value = 0.97697169189677
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 32.709374
$ws.Range("N9").Value = 98.12812199999999
$ws.Range("O9").Value = 0.6663257858061865
$ws.Range("P9").Value = 0.6663257858061865
$ws.Range("Q9").Value = 193.531749402044
$ws.Range("R9").Value = 1741.785744618396
$ws.Range("S9").Value = 0.6509814303135147
$ws.Range("T9").Value = 0.6509814303135149

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Edn3"
$ws.Range("C10").Value = "Ednra"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.916706
$ws.Range("H10").Value = 17.750118
$ws.Range("I10").Value = 0.9769716918967699
$ws.Range("J10").Value = 0.97697169189677
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.20978866666667
$ws.Range("N10").Value = 39.629366
$ws.Range("O10").Value = 0.2690978682028682
$ws.Range("P10").Value = 0.2690978682028682
$ws.Range("Q10").Value = 78.15843586279867
$ws.Range("R10").Value = 703.4259227651879
$ws.Range("S10").Value = 0.2629009995839701
$ws.Range("T10").Value = 0.2629009995839702
